# ROYALTY SERVICE REPORT - update to week # 15 period, add a daily receipt
# line (Carol Henry / Edwin Cabrera / Robert Melgoza, $1000 cash) plus the
# corresponding Daily Totals row and the expanded r/s-fee breakdown text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header block: franchisee / period / week number
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Franchisee: Amto.Robert (Arthur Murray Thousand Oaks, Asansol)"
$ws.Range("H2").Value = "(04/06/2025 - 04/12/2025)"
$ws.Range("K2").Value = "Week # 15"

# ---------------------------------------------------------------------
# 2. Insert two new rows above the old row 6 (the refunds disclaimer),
#    making room for a receipts detail line + a "Daily Totals" line.
# ---------------------------------------------------------------------
$ws.Rows("5:6").Insert()

# --- Row 5: the day's single receipt line -----------------------------
$ws.Cells.Item(5, 1).Value = 100000024
$ws.Cells.Item(5, 2).Value = "'04/09/2025"
$ws.Cells.Item(5, 3).Value = "Carol Henry"
$ws.Cells.Item(5, 4).Value = "Cash"
$ws.Cells.Item(5, 5).Value = "Edwin Cabrera"
$ws.Cells.Item(5, 6).Value = "Robert Melgoza"
$ws.Cells.Item(5, 7).Value = "37/REN"
$ws.Cells.Item(5, 8).Value = "0 / `$1000.00"
$ws.Cells.Item(5, 9).Value = 1000
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 1000
$ws.Cells.Item(5, 13).Value = 1000
$ws.Range("A5:M5").VerticalAlignment = -4108

# --- Row 6: "Daily Totals" summary line --------------------------------
$ws.Cells.Item(6, 1).Value = "Daily Totals"
$ws.Cells.Item(6, 8).Value = 1000
$ws.Cells.Item(6, 9).Value = 1000
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 1000
$ws.Cells.Item(6, 13).Value = 1000
$ws.Range("A6:G6").Merge()
$ws.Range("A6:M6").Font.Bold = $true

# ---------------------------------------------------------------------
# 3. "Total receipts" row (old row 11, now row 13): fill in the $1000
#    studio-receipts total that ripples through the r/s-fee section.
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 3).Value = 1000
$ws.Cells.Item(13, 8).Value = "'1,000.00"
$ws.Cells.Item(13, 13).Value = "'1,000.00"

# ---------------------------------------------------------------------
# 4. "Total subject to r/s fee" row (old row 13, now row 15): replace the
#    generic placeholder with the real fee computation breakdown.
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 5).Value = "Arthur Murray Thousand Oaks - "
$ws.Cells.Item(15, 8).Value = "`$1,000.00 ---------- 1000.00"
$ws.Cells.Item(15, 10).Value = "X 7 % - "
$ws.Cells.Item(15, 12).Value = "70.00 - "

# ---------------------------------------------------------------------
# 5. Selection cursor follows the last cell of the report, as in the
#    original workbook (always the bottom-right of the final row).
# ---------------------------------------------------------------------
$ws.Range("L15").Select()
